$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: turn the old lone date cell (A3) into a full new user record,
# matching the shared-string reuse/order Excel produces when a cell's
# old string reference is dropped and new strings are entered in order.
$ws.Range("A3").Value = 4
$ws.Range("B3").Value = "Кудлай"
$ws.Range("C3").Value = "Полина"
$ws.Range("D3").Value = "Александровна"
$ws.Range("E3").Value = "Иноватика"
$ws.Range("F3").Value = "pol@ina.com"

# G3 and K3 are purely-numeric-looking text ("789654", "555") that must be
# stored as shared-string text (no numeric value, no cell style change).
# Entering them directly would auto-infer a number. Route the literal text
# through a helper cell's formula result and paste-special (values only),
# which preserves the string type without adding any NumberFormat/style.
$ws.Range("Z1").Formula = "=""789654"""
$ws.Range("Z1").Copy()
$ws.Range("G3").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()

$ws.Range("H3").Value = "15.11.1999"

$ws.Range("Z1").Formula = "=""555"""
$ws.Range("Z1").Copy()
$ws.Range("K3").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()

# New row 4: a new "last visit" style timestamp string.
$ws.Range("A4").Value = "2019-08-03 00:00:00"
